$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ScenarioResults": drop the "Scenario number" column, drop the last
# data row, and replace the remaining data with the new scenario run results.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ScenarioResults")
[void]$ws1.Select()

# Remove column A ("Scenario number") - shifts B:F left into A:E.
[void]$ws1.Columns.Item(1).Delete()

# Remove the now-trailing row 6, leaving rows 1-5.
[void]$ws1.Rows.Item(6).Delete()

# Header row
$ws1.Range("A1").Value = "Url"
$ws1.Range("B1").Value = "Get requests"
$ws1.Range("C1").Value = "Post requests"
$ws1.Range("D1").Value = "Finded by id"
$ws1.Range("E1").Value = "Finded by tag"

# Row 2
$ws1.Range("A2").Value = "https://text.ru/"
$ws1.Range("B2").Value = 200
$ws1.Range("C2").Value = 200
$ws1.Range("D2").Value = "The element IS on the page"
$ws1.Range("E2").Value = "The operation was not requested"

# Row 3
$ws1.Range("A3").Value = "https://ecostyle.ua/pay/login.php?account=&phone="
$ws1.Range("B3").Value = "The operation was not requested"
$ws1.Range("C3").Value = "The operation was not requested"
$ws1.Range("D3").Value = "The element IS on the page"
$ws1.Range("E3").Value = "Tag IS here"

# Row 4
$ws1.Range("A4").Value = "https://text.ru/"
$ws1.Range("B4").Value = 200
$ws1.Range("C4").Value = 200
$ws1.Range("D4").Value = "The element IS on the page"
$ws1.Range("E4").Value = "The operation was not requested"

# Row 5
$ws1.Range("A5").Value = "https://ecostyle.ua/pay/login.php?account=&phone="
$ws1.Range("B5").Value = "The operation was not requested"
$ws1.Range("C5").Value = "The operation was not requested"
$ws1.Range("D5").Value = "The element IS on the page"
$ws1.Range("E5").Value = "Tag IS here"

[void]$ws1.Range("I9").Select()

# ---------------------------------------------------------------------------
# Sheet "ElementById": append the new scenario test rows 2-9.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ElementById")
[void]$ws4.Select()

$ws4.Range("A2").Value = "contact-form-7-js-extra"
$ws4.Range("B2").Value = "https://ecostyle.ua/pay/login.php?account=&phone="
$ws4.Range("C2").Value = "The element ISN'T on the page"

$ws4.Range("A3").Value = "contact-form-7-js-extra"
$ws4.Range("B3").Value = "https://ecostyle.ua/pay/login.php?account=&phone="
$ws4.Range("C3").Value = "The element IS on the page"

$ws4.Range("A4").Value = "contact-form-7-js-extr"
$ws4.Range("B4").Value = "https://ecostyle.ua/pay/login.php?account=&phone="
$ws4.Range("C4").Value = "The element ISN'T on the page"

$ws4.Range("A5").Value = "flags"
$ws4.Range("B5").Value = "https://ecostyle.ua/pay/login.php?account=&phone="
$ws4.Range("C5").Value = "The element IS on the page"

$ws4.Range("A6").Value = "hfcr"
$ws4.Range("B6").Value = "https://www.google.com/"
$ws4.Range("C6").Value = "The element ISN'T on the page"

$ws4.Range("A7").Value = "master-menu"
$ws4.Range("B7").Value = "https://text.ru/"
$ws4.Range("C7").Value = "The element IS on the page"

$ws4.Range("A8").Value = "master-menu"
$ws4.Range("B8").Value = "https://text.ru/"
$ws4.Range("C8").Value = "The element IS on the page"

$ws4.Range("A9").Value = "fruit-menu"
$ws4.Range("B9").Value = "https://text.ru/"
$ws4.Range("C9").Value = "The element IS on the page"

[void]$ws1.Select()
